$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.012165647470493
$ws.Range("D2").Value = 1.014811027772548
$ws.Range("E2").Value = 1.014247446349982
$ws.Range("F2").Value = 1.014780947268809
$ws.Range("I2").Value = 1.023038559704333
$ws.Range("J2").Value = 1.0174099303484
$ws.Range("K2").Value = 1.017668062636548
$ws.Range("L2").Value = 1.0171061698744
$ws.Range("M2").Value = 1.017638072211476
$ws.Range("N2").Value = 1.018854769148045

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.013911306951471
$ws.Range("D3").Value = 1.016467873830169
$ws.Range("E3").Value = 1.015744682139577
$ws.Range("F3").Value = 1.017689380109776
$ws.Range("I3").Value = 1.023341304150576
$ws.Range("J3").Value = 1.018784939528664
$ws.Range("K3").Value = 1.019127572552041
$ws.Range("L3").Value = 1.018406387168591
$ws.Range("M3").Value = 1.020345697098052
$ws.Range("N3").Value = 1.02023173099906

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.015028022685198
$ws.Range("D4").Value = 1.017527532453945
$ws.Range("E4").Value = 1.016702796726788
$ws.Range("F4").Value = 1.019518802689477
$ws.Range("I4").Value = 1.023523697959828
$ws.Range("J4").Value = 1.01966169489643
$ws.Range("K4").Value = 1.020059127607024
$ws.Range("L4").Value = 1.019236562013047
$ws.Range("M4").Value = 1.022045173828678
$ws.Range("N4").Value = 1.021109731460014

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.015494476271415
$ws.Range("D5").Value = 1.017970095236988
$ws.Range("E5").Value = 1.017103078186925
$ws.Range("F5").Value = 1.020275525039779
$ws.Range("I5").Value = 1.023597172999848
$ws.Range("J5").Value = 1.02002723344293
$ws.Range("K5").Value = 1.020447734246874
$ws.Range("L5").Value = 1.01958294663798
$ws.Range("M5").Value = 1.022747255763258
$ws.Range("N5").Value = 1.021475789113178

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.015572620861958
$ws.Range("D6").Value = 1.018044234017085
$ws.Range("E6").Value = 1.017170141451349
$ws.Range("F6").Value = 1.020401862657957
$ws.Range("I6").Value = 1.023609322736973
$ws.Range("J6").Value = 1.020088431693858
$ws.Range("K6").Value = 1.020512807540291
$ws.Range("L6").Value = 1.01964095382505
$ws.Range("M6").Value = 1.02286441810649
$ws.Range("N6").Value = 1.021537074272642

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.015034267223736
$ws.Range("D7").Value = 1.01753345738836
$ws.Range("E7").Value = 1.016708155103394
$ws.Range("F7").Value = 1.019528962375987
$ws.Range("I7").Value = 1.023524692286537
$ws.Range("J7").Value = 1.019666591151257
$ws.Range("K7").Value = 1.020064331978332
$ws.Range("L7").Value = 1.019241200656053
$ws.Range("M7").Value = 1.022054603462768
$ws.Range("N7").Value = 1.021114634668085

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.01275829815547
$ws.Range("D8").Value = 1.015373576949917
$ws.Range("E8").Value = 1.014755690732203
$ws.Range("F8").Value = 1.015774871811594
$ws.Range("I8").Value = 1.023143687643323
$ws.Range("J8").Value = 1.017877340333717
$ws.Range("K8").Value = 1.018164005428259
$ws.Range("L8").Value = 1.017547924114381
$ws.Range("M8").Value = 1.018564129348913
$ws.Range("N8").Value = 1.019322842909153

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.008646478980307
$ws.Range("D9").Value = 1.011469580297998
$ws.Range("E9").Value = 1.011230828879749
$ws.Range("F9").Value = 1.008747308339223
$ws.Range("I9").Value = 1.022367454386645
$ws.Range("J9").Value = 1.014622491928608
$ws.Range("K9").Value = 1.014714312314066
$ws.Range("L9").Value = 1.014476382028318
$ws.Range("M9").Value = 1.01200142837174
$ws.Range("N9").Value = 1.016063372245973

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.005832833329844
$ws.Range("D10").Value = 1.008796833167907
$ws.Range("E10").Value = 1.008820543878717
$ws.Range("F10").Value = 1.003769822999736
$ws.Range("I10").Value = 1.021777237823404
$ws.Range("J10").Value = 1.012380079595129
$ws.Range("K10").Value = 1.012342526683091
$ws.Range("L10").Value = 1.012366147546186
$ws.Range("M10").Value = 1.00733467180681
$ws.Range("N10").Value = 1.013817775429773

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.004596263617348
$ws.Range("D11").Value = 1.007621867930976
$ws.Range("E11").Value = 1.007761659066393
$ws.Range("F11").Value = 1.001541439819525
$ws.Range("I11").Value = 1.021503858276539
$ws.Range("J11").Value = 1.011390929186265
$ws.Range("K11").Value = 1.011297463688668
$ws.Range("L11").Value = 1.011436703752643
$ws.Range("M11").Value = 1.00524117088926
$ws.Range("N11").Value = 1.012827220313871

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.00413410963464
$ws.Range("D12").Value = 1.007182689123523
$ws.Range("E12").Value = 1.007365975039413
$ws.Range("F12").Value = 1.000702393736418
$ws.Range("I12").Value = 1.021399584338426
$ws.Range("J12").Value = 1.011020696270859
$ws.Range("K12").Value = 1.010906476938684
$ws.Range("L12").Value = 1.011089029886262
$ws.Range("M12").Value = 1.004452288751216
$ws.Range("N12").Value = 1.012456461625258

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.004233373220064
$ws.Range("D13").Value = 1.007277020215836
$ws.Range("E13").Value = 1.007450959079233
$ws.Range("F13").Value = 1.000882890230933
$ws.Range("I13").Value = 1.021422075816191
$ws.Range("J13").Value = 1.011100241509143
$ws.Range("K13").Value = 1.010990473309029
$ws.Range("L13").Value = 1.01116371870927
$ws.Range("M13").Value = 1.004622021867884
$ws.Range("N13").Value = 1.012536119826905

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.004558120274623
$ws.Range("D14").Value = 1.007585621816855
$ws.Range("E14").Value = 1.00772900047544
$ws.Range("F14").Value = 1.001472317108231
$ws.Range("I14").Value = 1.021495294941363
$ws.Range("J14").Value = 1.011360383643868
$ws.Range("K14").Value = 1.011265202338926
$ws.Range("L14").Value = 1.01140801510943
$ws.Range("M14").Value = 1.005176193393203
$ws.Range("N14").Value = 1.0127966313933

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.004757828754639
$ws.Range("D15").Value = 1.007775394909173
$ws.Range("E15").Value = 1.007899994786214
$ws.Range("F15").Value = 1.00183397108711
$ws.Range("I15").Value = 1.021540044453077
$ws.Range("J15").Value = 1.0115202895921
$ws.Range("K15").Value = 1.011434097642922
$ws.Range("L15").Value = 1.011558208816757
$ws.Range("M15").Value = 1.005516134070689
$ws.Range("N15").Value = 1.01295676442632

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.005914507962533
$ws.Range("D16").Value = 1.008874432111595
$ws.Range("E16").Value = 1.008890491195339
$ws.Range("F16").Value = 1.00391614458611
$ws.Range("I16").Value = 1.021795001453412
$ws.Range("J16").Value = 1.012445335987758
$ws.Range("K16").Value = 1.012411495937899
$ws.Range("L16").Value = 1.012427494330877
$ws.Range("M16").Value = 1.007472049796538
$ws.Range("N16").Value = 1.013883124493964

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.006635112148766
$ws.Range("D17").Value = 1.009559040388045
$ws.Range("E17").Value = 1.009507674122961
$ws.Range("F17").Value = 1.005202431385152
$ws.Range("I17").Value = 1.021950124661035
$ws.Range("J17").Value = 1.013020667580809
$ws.Range("K17").Value = 1.013019694094447
$ws.Range("L17").Value = 1.012968517364111
$ws.Range("M17").Value = 1.008679236589768
$ws.Range("N17").Value = 1.014459273123868

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.007053675044783
$ws.Range("D18").Value = 1.009956664615718
$ws.Range("E18").Value = 1.009866204153707
$ws.Range("F18").Value = 1.005945672109357
$ws.Range("I18").Value = 1.02203889058719
$ws.Range("J18").Value = 1.013354501791362
$ws.Range("K18").Value = 1.01337270923263
$ws.Range("L18").Value = 1.013282577853805
$ws.Range("M18").Value = 1.009376371041654
$ws.Range("N18").Value = 1.014793581417292

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.007196099603628
$ws.Range("D19").Value = 1.01009195922532
$ws.Range("E19").Value = 1.009988207972254
$ws.Range("F19").Value = 1.006197915851226
$ws.Range("I19").Value = 1.022068868140796
$ws.Range("J19").Value = 1.013468036905287
$ws.Range("K19").Value = 1.013492786357886
$ws.Range("L19").Value = 1.013389410700543
$ws.Range("M19").Value = 1.009612898699405
$ws.Range("N19").Value = 1.014907277764101

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.006557980171656
$ws.Range("D20").Value = 1.009485764498548
$ws.Range("E20").Value = 1.009441608061717
$ws.Range("F20").Value = 1.005065154392007
$ws.Range("I20").Value = 1.021933659147563
$ws.Range("J20").Value = 1.012959121264826
$ws.Range("K20").Value = 1.012954620440218
$ws.Range("L20").Value = 1.012910627357588
$ws.Range("M20").Value = 1.008550442940776
$ws.Range("N20").Value = 1.014397639405058

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.004462569476825
$ws.Range("D21").Value = 1.007494822872056
$ws.Range("E21").Value = 1.007647190233651
$ws.Range("F21").Value = 1.001299061091617
$ws.Range("I21").Value = 1.021473809499484
$ws.Range("J21").Value = 1.011283856809561
$ws.Range("K21").Value = 1.011184379636942
$ws.Range("L21").Value = 1.011336143839262
$ws.Range("M21").Value = 1.005013317366617
$ws.Range("N21").Value = 1.012719995882111

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.0031286359874
$ws.Range("D22").Value = 1.006227110253995
$ws.Range("E22").Value = 1.006505229543508
$ws.Range("F22").Value = 0.9988654566840138
$ws.Range("I22").Value = 1.021168867884034
$ws.Range("J22").Value = 1.010214197152283
$ws.Range("K22").Value = 1.010055088217417
$ws.Range("L22").Value = 1.010332060925878
$ws.Range("M22").Value = 1.002724043203645
$ws.Range("N22").Value = 1.011648817185427

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.003837373803135
$ws.Range("D23").Value = 1.006900691127567
$ws.Range("E23").Value = 1.007111935228735
$ws.Range("F23").Value = 1.000161905883748
$ws.Range("I23").Value = 1.021332041000067
$ws.Range("J23").Value = 1.010782825197304
$ws.Range("K23").Value = 1.010655320766712
$ws.Range("L23").Value = 1.010865712380482
$ws.Range("M23").Value = 1.003943940924817
$ws.Range("N23").Value = 1.012218252747495

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.006592838203988
$ws.Range("D24").Value = 1.009518879954498
$ws.Range("E24").Value = 1.009471464984715
$ws.Range("F24").Value = 1.005127205652192
$ws.Range("I24").Value = 1.021941104499826
$ws.Range("J24").Value = 1.012986936788503
$ws.Range("K24").Value = 1.012984029783959
$ws.Range("L24").Value = 1.012936790017987
$ws.Range("M24").Value = 1.008608660846629
$ws.Range("N24").Value = 1.014425494429968

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.009721914044853
$ws.Range("D25").Value = 1.012490885948267
$ws.Range("E25").Value = 1.01215244972749
$ws.Range("F25").Value = 1.010614466441117
$ws.Range("I25").Value = 1.022580772358041
$ws.Range("J25").Value = 1.015476422491119
$ws.Range("K25").Value = 1.015618517275667
$ws.Range("L25").Value = 1.015281199431104
$ws.Range("M25").Value = 1.013748308825956
$ws.Range("N25").Value = 1.01691851548783
